# UITMBER.docx edit script
#
# 1) Insert a new "DDD - Data Driven Development" Heading-3 paragraph plus a
#    blank Heading-3 paragraph right before the existing "Zalozenia
#    funkcjonalne" heading.
# 2) Move the <w:lastRenderedPageBreak/> marker from the run that currently
#    carries it ("Saldo/dzienny przychod") to the run right before it in the
#    list ("Zdjecie pojazdu,").

$d = $word.ActiveDocument

function Find-ParagraphIndexByText($doc, [string]$needle) {
    # Locate $needle anywhere in the story, then resolve it back to the
    # index of the Paragraphs collection it lives in (Find's own hit-range
    # doesn't carry a usable live Paragraphs collection in this host, so we
    # map the match's Start/End back onto $doc.Paragraphs ourselves).
    $hit = $doc.Content
    $ok = $hit.Find.Execute($needle, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        return -1
    }
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $cand = $doc.Paragraphs.Item($i)
        if ($cand.Range.Start -le $hit.Start -and $cand.Range.End -ge $hit.End) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) "DDD - Data Driven Development" + blank Heading-3 paragraph, inserted
#    immediately before "Zalozenia funkcjonalne".
# ---------------------------------------------------------------------

$zalozeniaIdx = Find-ParagraphIndexByText $d "Założenia funkcjonalne"
$zalozeniaPara = $d.Paragraphs.Item($zalozeniaIdx)
$zalozeniaRange = $zalozeniaPara.Range

# Two fresh paragraphs ahead of it; InsertParagraphBefore copies the
# paragraph formatting (Heading 3 / Nagwek3 style) of the paragraph it is
# attached to, which is exactly what we want here.
$zalozeniaRange.InsertParagraphBefore()
$zalozeniaRange.InsertParagraphBefore()

$dddIdx = Find-ParagraphIndexByText $d "Założenia funkcjonalne"
$dddIdx = $dddIdx - 2
$dddPara = $d.Paragraphs.Item($dddIdx)
$dddPara.Range.Text = "DDD – Data Driven Development"
# (the paragraph right after $dddPara, before "Zalozenia funkcjonalne", is
# left as the required blank Heading-3 paragraph)

# ---------------------------------------------------------------------
# 2) Relocate <w:lastRenderedPageBreak/> from "Saldo/dzienny przychod" to
#    "Zdjecie pojazdu,".
# ---------------------------------------------------------------------

$zdjecieIdx = Find-ParagraphIndexByText $d "Zdjęcie pojazdu,"
$zdjeciePara = $d.Paragraphs.Item($zdjecieIdx)
$zdjecieRange = $zdjeciePara.Range

# Splice in a brand-new paragraph ahead of the old one and fill it via raw
# OOXML so the otherwise COM-unreachable <w:lastRenderedPageBreak/> run
# child can be set explicitly.
$zdjecieRange.InsertParagraphBefore()
$newZdjecieIdx = $zdjecieIdx
$newZdjeciePara = $d.Paragraphs.Item($newZdjecieIdx)
$newZdjecieRange = $newZdjeciePara.Range

$zdjeciePkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Zdjęcie pojazdu,</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$newZdjecieRange.InsertXML($zdjeciePkg)

# The old (page-break-less) "Zdjecie pojazdu," paragraph got pushed one
# slot down; remove it now that its replacement is in place.
$oldZdjecieIdx = $newZdjecieIdx + 1
$oldZdjeciePara = $d.Paragraphs.Item($oldZdjecieIdx)
$oldZdjeciePara.Range.Delete()

# Strip the marker from its old home: re-assigning the run's text forces a
# fresh run without the lastRenderedPageBreak child, while leaving the
# paragraph's own formatting (style/numPr) untouched.
$saldoIdx = Find-ParagraphIndexByText $d "Saldo/dzienny przychód"
$saldoPara = $d.Paragraphs.Item($saldoIdx)
$saldoPara.Range.Text = "Saldo/dzienny przychód"

Write-Output "edit complete"
